$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column G: "utilities" -> "basement"
$ws.Range("G1").Value = "basement"

# Update values in column G: "landlord"/"tenant" -> "finished"/"unfinished"
$ws.Range("G2").Value = "finished"
$ws.Range("G3").Value = "finished"
$ws.Range("G4").Value = "unfinished"

# Set column G ("basement") width to match the best-fit sizing Excel computed
# for its new content ("basement" / "finished" / "unfinished").
$ws.Columns.Item(7).ColumnWidth = 8.5

# Update selected cell
$ws.Range("H10").Select()
